{"js": "const pairs = [\n  [\"[Music]\", \"[\u1219\u12da\u1243]\"],\n  [\"okay so the puzzles I'm going to\", \"\u12a5\u123a \u12e8\u121d\u1290\u130d\u122b\u127d\u1201 \u12a5\u1290\u1246\u1245\u120d\u123d\"],\n  [\"challenge you with are two basic\", \"\u1201\u1208\u1275 \u1218\u1230\u1228\u1273\u12ca\"],\n  [\"versions of a more complicated puzzle\", \"\u1235\u122a\u1276\u127d \u12e8\u1206\u1291\u1275\u1295 \u1260\u1323\u121d \u12cd\u1235\u1265\u1235\u1265 \u12e8\u1206\u1290 \u12a5\u1295\u1246\u1245\u120d\u123d\"],\n  [\"known as the ants puzzle, which I'm\", \"\u12e8\u1309\u1295\u12f3\u1296\u127d \u12a5\u1290\u1246\u1245\u120d\u123d \u1235\u1208\u121a\u1263\u1208\u12cd\u1363\"],\n  [\"probably going to discuss in a different\", \"\u121d\u1293\u120d\u1263\u1275 \u12e8\u121b\u12c8\u122b\u1260\u1275 \u12ed\u1206\u1293\u120d \u1260\u120c\u120b\"],\n  [\"video. Let me just finish writing down\", \"\u126a\u12f5\u12ee\u1361\u1361 \u12a0\u1201\u1295 \u1345\u134c \u120d\u1328\u122d\u1235\"],\n  [\"the title and, well, I can even draw a\", \"\u122d\u12d5\u1231\u1295\u1363 \u1218\u120d\u12ab\u121d\u1363 \u120d\u1235\u120d\u120b\u127d\u1201\u121d \u12a5\u127d\u120b\u1208\u1201\"],\n  [\"little ant right here. okay, let's get\", \"\u1275\u1295\u123d \u1309\u1295\u12f3\u1295 \u12a5\u12da\u1205\u1361\u1361 \u12a5\u123a\"],\n  [\"started! As I said I'm going to discuss\", \"\u12a5\u1295\u1300\u121d\u122d\u1361\u1361  \u12a5\u1295\u12f3\u120d\u12b3\u127d\u1201 \u12a8\u121d\u1290\u130d\u122b\u127d\u1201\"],\n  [\"two puzzles in the first puzzle there\", \"\u1201\u1208\u1275 \u12a5\u1290\u1246\u1245\u120d\u123e\u127d \u1260\u1218\u1300\u1218\u122a\u12eb\u12cd \u12a5\u1295\u1246\u1245\u120d\u123d\"],\n  [\"are two ants on a very high stool: a sort\", \"\u1201\u1208\u1275 \u1309\u1295\u12f3\u1296\u127d \u1260\u1323\u121d \u12a8\u134d \u12ab\u1208 \u1260\u122d\u1329\u121b \u120b\u12ed\u1361- \u1270\u122b\u122b\"],\n  [\"of Mountain, flat on the top with two\", \"\u12a0\u12ed\u1290\u1275\u1363 \u1320\u134d\u1323\u134b \u12a0\u1293\u1271 \u120b\u12ed \u12a8\u1201\u1208\u1275\"],\n  [\"steep cliffs to both the sides. The flat\", \"\u1241\u120d\u1241\u1208\u1273\u121b \u1320\u122d\u12dd \u130b\u122d \u1260\u1201\u1208\u1275 \u1260\u12a9\u120d\u1361\u1361 \u1320\u134d\u1323\u134b\u12cd\"],\n  [\"peak is one meter wide the two ants move\", \"\u132b\u134d \u12a0\u1295\u12f5 \u121c\u1275\u122d \u1235\u134b\u1275 \u12a0\u1208\u12cd \u1201\u1208\u1271 \u1309\u1295\u12f3\u1296\u127d \u1270\u1295\u1240\u1233\u1240\u1231\"],\n  [\"with a velocity, let's call it V, which is\", \"\u1260\u1206\u1290 \u134d\u1325\u1290\u1275\u1363 \u134d \u12a5\u1295\u1260\u1208\u12cd\u1363 \u12ed\u1205\"],\n  [\"the same for both of them and that is\", \"\u1208\u1201\u1208\u1271\u121d \u1270\u1218\u1233\u1233\u12ed \u1290\u12cd \u12a5\u1293 \u12ed\u1205\u121d\"],\n  [\"equal to one centimeter per second. You\", \"\u12a0\u1295\u12f5 \u1234\u1295\u1272\u121c\u1275\u122d \u1260\u1230\u12ae\u1295\u12f5 \u1290\u12cd\u1361\u1361 \u12a5\u1293\u1295\u1270\"],\n  [\"can decide the direction towards each\", \"\u1218\u12c8\u1230\u1295 \u1275\u127d\u120b\u120b\u127d\u1201 \u12e8\u121a\u1313\u12d9\u1260\u1275\u1295 \u12a0\u1245\u1323\u132b \u12e8\u12a5\u12eb\u1295\u12f3\u1295\u12f1\u1295\"],\n  [\"ant moves if it is right or left and\", \"\u1309\u1295\u12f3\u1295 \u12a5\u1295\u1245\u1235\u1243\u1234 \u1240\u129d \u12a8\u1206\u1290 \u12c8\u12ed\u121d \u130d\u122b \u12a5\u1293\"],\n  [\"where exactly to place the two ants on the\", \"\u1201\u1208\u1271 \u1309\u1295\u12f3\u1296\u127d \u1260\u1275\u12ad\u12ad\u120d \u12e8\u1275 \u1218\u1206\u1295 \u12a5\u1295\u12f3\u1208\u1263\u1278\u12cd\"],\n  [\"top of the mountain. Your purpose is to\", \"\u1260\u1270\u122b\u122b\u12cd \u132b\u134d \u120b\u12ed\u1361\u1361 \u12e8\u1293\u1295\u1270 \u12a0\u120b\u121b\"],\n  [\"make the time the last ant takes before\", \"\u12e8\u1218\u1328\u1228\u123b\u12cd \u1309\u1295\u12f3\u1295 \u12a8\u1218\u12cd\u12f0\u1241 \u1260\u134a\u1275\"],\n  [\"falling the longest possible. Ants cannot\", \"\u12eb\u1208\u12cd\u1295 \u130a\u12dc \u1228\u1305\u121d \u121b\u12f5\u1228\u130d \u1290\u12cd\u1361\u1361 \u1309\u1295\u12f3\u1296\u1279\"],\n  [\"be still: they must move to the right or\", \"\u12a0\u12ed\u1246\u1219\u121d\u1361- \u1218\u1295\u1240\u1233\u1240\u1235 \u12a0\u1208\u1263\u1278\u12cd \u12c8\u12f0 \u130d\u122b \u12c8\u12ed\u121d\"],\n  [\"to the left but they must move and after\", \"\u12c8\u12f0 \u1240\u129d \u130d\u1295 \u1218\u1295\u1240\u1233\u1240\u1235 \u12a0\u1208\u1263\u1278\u12cd \u12a5\u1293\"],\n  [\"meeting each other they turn around and\", \"\u12a8\u1270\u1308\u1293\u1299 \u1260\u128b\u120b \u12ed\u12de\u1229\u1293\"],\n  [\"keep moving with the same but opposite\", \"\u1218\u1295\u1240\u1233\u1240\u1235 \u12ed\u1240\u1325\u120b\u1209 \u1260\u1270\u1218\u1233\u1233\u12ed \u130d\u1295 \u1260\u1270\u1243\u122b\u1292\"],\n  [\"velocity\", \"\u134d\u1325\u1290\u1275\"],\n  [\"[Music]\", \"[\u1219\u12da\u1243]\"],\n  [\"so again what are the precise positions\", \"\u1235\u1208\u12da\u1205 \u12a0\u1201\u1295\u121d \u12e8\u1275 \u1293\u1278\u12cd \u1275\u12ad\u12ad\u1208\u129b \u1266\u1273\u12ce\u1279\"],\n  [\"where I should place the two ants in\", \"\u1201\u1208\u1271\u1295 \u1309\u1295\u12f3\u1296\u127d \u12e8\u121b\u1235\u1240\u121d\u1325\u1260\u1275\"],\n  [\"order to get the longest time before the\", \"\u1228\u1305\u121d \u130a\u12dc \u1208\u121b\u1308\u1298\u1275\"],\n  [\"last ant falls? The second puzzle is\", \"\u12e8\u1218\u1328\u1228\u123b\u12cd \u1309\u1295\u12f3\u1295 \u12a8\u1218\u12cd\u12f0\u1241 \u1260\u134a\u1275? \u1201\u1208\u1270\u129b\u12cd \u12a5\u1295\u1246\u1245\u120d\u123d\"],\n  [\"basically the same but now we have three\", \"\u1260\u1218\u1230\u1228\u1273\u12ca\u1290\u1275 \u1270\u1218\u1233\u1233\u12ed \u1290\u12cd \u130d\u1295 \u12a0\u1201\u1295 \u1236\u1235\u1275\"],\n  [\"ants instead of two.\", \"\u1309\u1295\u12f3\u1296\u127d \u12a0\u1209 \u1260\u1201\u1208\u1271 \u1266\u1273\"],\n  [\"As before the ants velocity is one\", \"\u12a5\u1295\u12f0\u1260\u134a\u1271 \u12e8\u1309\u1295\u12f3\u1296\u1279 \u134d\u1325\u1290\u1275 \u12a0\u1295\u12f5\"],\n  [\"centimeter per second, every ant turns\", \"\u1234\u1295\u1272\u121c\u1275\u122d \u1260\u1230\u12a8\u1295\u12f5\u1363 \u12a5\u1295\u12f3\u1295\u12f1 \u1309\u1295\u12f3\u1295 \u12ed\u1218\u1208\u1233\u120d\"],\n  [\"around after meeting another ant and\", \"\u12a8\u120c\u120b \u1309\u1295\u12f3\u1295 \u130b\u122d \u12a8\u1270\u1308\u1293\u1298 \u1260\u128b\u120b \u12a5\u1293\"],\n  [\"the peak is one meter wide. So, what are\", \"\u12a8\u134d\u1273\u12cd \u12a0\u1295\u12f5 \u121c\u1275\u122d \u1235\u134b\u1275 \u12a0\u1208\u12cd\u1361\u1361 \u1235\u1208\u12da\u1205\u1363 \u12e8\u1275 \u1293\u1278\u12cd\"],\n  [\"now the precise positions\", \"\u12a0\u1201\u1295 \u1260\u1275\u12ad\u12ad\u120d \u1266\u1273\u12ce\u1279\"],\n  [\"I should place the three ants in order\", \"\u1236\u1235\u1271\u1295 \u1309\u1295\u12f3\u1296\u127d \u121b\u1235\u1240\u1218\u1325 \u12eb\u1208\u1265\u129d\"],\n  [\"to get the longest time before the last\", \"\u1228\u1305\u1219\u1295 \u130a\u12dc \u1208\u121b\u1308\u1298\u1275 \u12e8\u1218\u1328\u1228\u123b\u12cd\"],\n  [\"ant falls down? I hope you enjoyed this\", \"\u1309\u1295\u12f3\u1295 \u12a8\u1218\u12cd\u12f0\u1241 \u1260\u134a\u1275?  \u12a5\u1295\u12f0\u1270\u12dd\u1293\u1293\u127d\u1201\u1260\u1275 \u1270\u1235\u134b \u12a0\u1228\u130b\u1208\u1201 \u1260\u12da\u1205\"],\n  [\"video do your best and good luck\", \"\u126a\u12f5\u12ee\u1361\u1361 \u12e8\u1270\u127b\u120b\u127d\u1201\u1295 \u12a0\u12f5\u122d\u1309 \u1218\u120d\u12ab\u121d \u12d5\u12f5\u120d\u1361\u1361\"],\n];\n\nconst paras = context.document.body.paragraphs;\nparas.load(\"items/text\");\nawait context.sync();\n\nfor (const [oldText, newText] of pairs) {\n  for (let i = 0; i < paras.items.length; i++) {\n    if (paras.items[i].text === oldText) {\n      paras.items[i].insertText(newText, \"Replace\");\n    }\n  }\n}\nawait context.sync();", "ps1": "$d = $word.ActiveDocument\n$pairs = @(\n    ,@('[Music]', '[\u1219\u12da\u1243]')\n    ,@('okay so the puzzles I''m going to', '\u12a5\u123a \u12e8\u121d\u1290\u130d\u122b\u127d\u1201 \u12a5\u1290\u1246\u1245\u120d\u123d')\n    ,@('challenge you with are two basic', '\u1201\u1208\u1275 \u1218\u1230\u1228\u1273\u12ca')\n    ,@('versions of a more complicated puzzle', '\u1235\u122a\u1276\u127d \u12e8\u1206\u1291\u1275\u1295 \u1260\u1323\u121d \u12cd\u1235\u1265\u1235\u1265 \u12e8\u1206\u1290 \u12a5\u1295\u1246\u1245\u120d\u123d')\n    ,@('known as the ants puzzle, which I''m', '\u12e8\u1309\u1295\u12f3\u1296\u127d \u12a5\u1290\u1246\u1245\u120d\u123d \u1235\u1208\u121a\u1263\u1208\u12cd\u1363')\n    ,@('probably going to discuss in a different', '\u121d\u1293\u120d\u1263\u1275 \u12e8\u121b\u12c8\u122b\u1260\u1275 \u12ed\u1206\u1293\u120d \u1260\u120c\u120b')\n    ,@('video. Let me just finish writing down', '\u126a\u12f5\u12ee\u1361\u1361 \u12a0\u1201\u1295 \u1345\u134c \u120d\u1328\u122d\u1235')\n    ,@('the title and, well, I can even draw a', '\u122d\u12d5\u1231\u1295\u1363 \u1218\u120d\u12ab\u121d\u1363 \u120d\u1235\u120d\u120b\u127d\u1201\u121d \u12a5\u127d\u120b\u1208\u1201')\n    ,@('little ant right here. okay, let''s get', '\u1275\u1295\u123d \u1309\u1295\u12f3\u1295 \u12a5\u12da\u1205\u1361\u1361 \u12a5\u123a')\n    ,@('started! As I said I''m going to discuss', '\u12a5\u1295\u1300\u121d\u122d\u1361\u1361  \u12a5\u1295\u12f3\u120d\u12b3\u127d\u1201 \u12a8\u121d\u1290\u130d\u122b\u127d\u1201')\n    ,@('two puzzles in the first puzzle there', '\u1201\u1208\u1275 \u12a5\u1290\u1246\u1245\u120d\u123e\u127d \u1260\u1218\u1300\u1218\u122a\u12eb\u12cd \u12a5\u1295\u1246\u1245\u120d\u123d')\n    ,@('are two ants on a very high stool: a sort', '\u1201\u1208\u1275 \u1309\u1295\u12f3\u1296\u127d \u1260\u1323\u121d \u12a8\u134d \u12ab\u1208 \u1260\u122d\u1329\u121b \u120b\u12ed\u1361- \u1270\u122b\u122b')\n    ,@('of Mountain, flat on the top with two', '\u12a0\u12ed\u1290\u1275\u1363 \u1320\u134d\u1323\u134b \u12a0\u1293\u1271 \u120b\u12ed \u12a8\u1201\u1208\u1275')\n    ,@('steep cliffs to both the sides. The flat', '\u1241\u120d\u1241\u1208\u1273\u121b \u1320\u122d\u12dd \u130b\u122d \u1260\u1201\u1208\u1275 \u1260\u12a9\u120d\u1361\u1361 \u1320\u134d\u1323\u134b\u12cd')\n    ,@('peak is one meter wide the two ants move', '\u132b\u134d \u12a0\u1295\u12f5 \u121c\u1275\u122d \u1235\u134b\u1275 \u12a0\u1208\u12cd \u1201\u1208\u1271 \u1309\u1295\u12f3\u1296\u127d \u1270\u1295\u1240\u1233\u1240\u1231')\n    ,@('with a velocity, let''s call it V, which is', '\u1260\u1206\u1290 \u134d\u1325\u1290\u1275\u1363 \u134d \u12a5\u1295\u1260\u1208\u12cd\u1363 \u12ed\u1205')\n    ,@('the same for both of them and that is', '\u1208\u1201\u1208\u1271\u121d \u1270\u1218\u1233\u1233\u12ed \u1290\u12cd \u12a5\u1293 \u12ed\u1205\u121d')\n    ,@('equal to one centimeter per second. You', '\u12a0\u1295\u12f5 \u1234\u1295\u1272\u121c\u1275\u122d \u1260\u1230\u12ae\u1295\u12f5 \u1290\u12cd\u1361\u1361 \u12a5\u1293\u1295\u1270')\n    ,@('can decide the direction towards each', '\u1218\u12c8\u1230\u1295 \u1275\u127d\u120b\u120b\u127d\u1201 \u12e8\u121a\u1313\u12d9\u1260\u1275\u1295 \u12a0\u1245\u1323\u132b \u12e8\u12a5\u12eb\u1295\u12f3\u1295\u12f1\u1295')\n    ,@('ant moves if it is right or left and', '\u1309\u1295\u12f3\u1295 \u12a5\u1295\u1245\u1235\u1243\u1234 \u1240\u129d \u12a8\u1206\u1290 \u12c8\u12ed\u121d \u130d\u122b \u12a5\u1293')\n    ,@('where exactly to place the two ants on the', '\u1201\u1208\u1271 \u1309\u1295\u12f3\u1296\u127d \u1260\u1275\u12ad\u12ad\u120d \u12e8\u1275 \u1218\u1206\u1295 \u12a5\u1295\u12f3\u1208\u1263\u1278\u12cd')\n    ,@('top of the mountain. Your purpose is to', '\u1260\u1270\u122b\u122b\u12cd \u132b\u134d \u120b\u12ed\u1361\u1361 \u12e8\u1293\u1295\u1270 \u12a0\u120b\u121b')\n    ,@('make the time the last ant takes before', '\u12e8\u1218\u1328\u1228\u123b\u12cd \u1309\u1295\u12f3\u1295 \u12a8\u1218\u12cd\u12f0\u1241 \u1260\u134a\u1275')\n    ,@('falling the longest possible. Ants cannot', '\u12eb\u1208\u12cd\u1295 \u130a\u12dc \u1228\u1305\u121d \u121b\u12f5\u1228\u130d \u1290\u12cd\u1361\u1361 \u1309\u1295\u12f3\u1296\u1279')\n    ,@('be still: they must move to the right or', '\u12a0\u12ed\u1246\u1219\u121d\u1361- \u1218\u1295\u1240\u1233\u1240\u1235 \u12a0\u1208\u1263\u1278\u12cd \u12c8\u12f0 \u130d\u122b \u12c8\u12ed\u121d')\n    ,@('to the left but they must move and after', '\u12c8\u12f0 \u1240\u129d \u130d\u1295 \u1218\u1295\u1240\u1233\u1240\u1235 \u12a0\u1208\u1263\u1278\u12cd \u12a5\u1293')\n    ,@('meeting each other they turn around and', '\u12a8\u1270\u1308\u1293\u1299 \u1260\u128b\u120b \u12ed\u12de\u1229\u1293')\n    ,@('keep moving with the same but opposite', '\u1218\u1295\u1240\u1233\u1240\u1235 \u12ed\u1240\u1325\u120b\u1209 \u1260\u1270\u1218\u1233\u1233\u12ed \u130d\u1295 \u1260\u1270\u1243\u122b\u1292')\n    ,@('velocity', '\u134d\u1325\u1290\u1275')\n    ,@('[Music]', '[\u1219\u12da\u1243]')\n    ,@('so again what are the precise positions', '\u1235\u1208\u12da\u1205 \u12a0\u1201\u1295\u121d \u12e8\u1275 \u1293\u1278\u12cd \u1275\u12ad\u12ad\u1208\u129b \u1266\u1273\u12ce\u1279')\n    ,@('where I should place the two ants in', '\u1201\u1208\u1271\u1295 \u1309\u1295\u12f3\u1296\u127d \u12e8\u121b\u1235\u1240\u121d\u1325\u1260\u1275')\n    ,@('order to get the longest time before the', '\u1228\u1305\u121d \u130a\u12dc \u1208\u121b\u1308\u1298\u1275')\n    ,@('last ant falls? The second puzzle is', '\u12e8\u1218\u1328\u1228\u123b\u12cd \u1309\u1295\u12f3\u1295 \u12a8\u1218\u12cd\u12f0\u1241 \u1260\u134a\u1275? \u1201\u1208\u1270\u129b\u12cd \u12a5\u1295\u1246\u1245\u120d\u123d')\n    ,@('basically the same but now we have three', '\u1260\u1218\u1230\u1228\u1273\u12ca\u1290\u1275 \u1270\u1218\u1233\u1233\u12ed \u1290\u12cd \u130d\u1295 \u12a0\u1201\u1295 \u1236\u1235\u1275')\n    ,@('ants instead of two.', '\u1309\u1295\u12f3\u1296\u127d \u12a0\u1209 \u1260\u1201\u1208\u1271 \u1266\u1273')\n    ,@('As before the ants velocity is one', '\u12a5\u1295\u12f0\u1260\u134a\u1271 \u12e8\u1309\u1295\u12f3\u1296\u1279 \u134d\u1325\u1290\u1275 \u12a0\u1295\u12f5')\n    ,@('centimeter per second, every ant turns', '\u1234\u1295\u1272\u121c\u1275\u122d \u1260\u1230\u12a8\u1295\u12f5\u1363 \u12a5\u1295\u12f3\u1295\u12f1 \u1309\u1295\u12f3\u1295 \u12ed\u1218\u1208\u1233\u120d')\n    ,@('around after meeting another ant and', '\u12a8\u120c\u120b \u1309\u1295\u12f3\u1295 \u130b\u122d \u12a8\u1270\u1308\u1293\u1298 \u1260\u128b\u120b \u12a5\u1293')\n    ,@('the peak is one meter wide. So, what are', '\u12a8\u134d\u1273\u12cd \u12a0\u1295\u12f5 \u121c\u1275\u122d \u1235\u134b\u1275 \u12a0\u1208\u12cd\u1361\u1361 \u1235\u1208\u12da\u1205\u1363 \u12e8\u1275 \u1293\u1278\u12cd')\n    ,@('now the precise positions', '\u12a0\u1201\u1295 \u1260\u1275\u12ad\u12ad\u120d \u1266\u1273\u12ce\u1279')\n    ,@('I should place the three ants in order', '\u1236\u1235\u1271\u1295 \u1309\u1295\u12f3\u1296\u127d \u121b\u1235\u1240\u1218\u1325 \u12eb\u1208\u1265\u129d')\n    ,@('to get the longest time before the last', '\u1228\u1305\u1219\u1295 \u130a\u12dc \u1208\u121b\u1308\u1298\u1275 \u12e8\u1218\u1328\u1228\u123b\u12cd')\n    ,@('ant falls down? I hope you enjoyed this', '\u1309\u1295\u12f3\u1295 \u12a8\u1218\u12cd\u12f0\u1241 \u1260\u134a\u1275?  \u12a5\u1295\u12f0\u1270\u12dd\u1293\u1293\u127d\u1201\u1260\u1275 \u1270\u1235\u134b \u12a0\u1228\u130b\u1208\u1201 \u1260\u12da\u1205')\n    ,@('video do your best and good luck', '\u126a\u12f5\u12ee\u1361\u1361 \u12e8\u1270\u127b\u120b\u127d\u1201\u1295 \u12a0\u12f5\u122d\u1309 \u1218\u120d\u12ab\u121d \u12d5\u12f5\u120d\u1361\u1361')\n)\n\nforeach ($pair in $pairs) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n    foreach ($p in $d.Paragraphs) {\n        $r = $p.Range\n        $txt = $r.Text\n        if ($txt.Length -gt 0 -and [int]$txt[$txt.Length-1] -eq 13) {\n            $txt = $txt.Substring(0, $txt.Length - 1)\n        }\n        if ($txt -eq $oldText) {\n            $r.Text = $newText\n        }\n    }\n}"}
